# =========================================================================
# Add "2022-Q4" worksheet (with fund-holdings detail) as the 2nd sheet,
# and add a corresponding summary row to the "总计" (total) sheet.
# =========================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (Total) sheet: insert a new row 2 for "2022-Q4"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Insert a new row at position 2, pushing existing data rows down by one.
$totalSheet.Rows.Item(2).Insert()

# The new row 2 lost formatting info from the insert; copy the style
# pattern from row 3 (the row that used to be row 2) onto it.
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

# Populate the new "2022-Q4" row.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 16
$totalSheet.Range("D2").Value = 1.56

# Re-sequence the running index in column A for every row that shifted
# down (old row N is now row N+1; the index itself also increments).
for ($r = 3; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right after "总计" (i.e. before
#    the worksheet that is currently in slot 2, "2022-Q2").
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# Re-resolve a fresh, name-based handle to the "2022-Q2" sheet (template)
# for copying styles -- positional handles can go stale once the sheet
# collection is mutated.
$templateSheet = $wb.Worksheets.Item("2022-Q2")

# Copy the header + first six data rows' style pattern (bold/bordered
# header row, bold/bordered column A, plain data cells) onto the new
# sheet's first seven rows.
$templateSheet.Range("A1:H7").Copy()
$newSheet.Range("A1:H7").PasteSpecial(-4122)

# Extend the "data row" style pattern (bold/bordered column A, plain
# data cells) down through row 17 (16 data rows total).
$templateSheet.Range("A2:H2").Copy()
$newSheet.Range("A8:H17").PasteSpecial(-4122)

# ---- Header row ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# Columns B (fund code), D, E, F, G (decimal-looking text values) must
# be stored as TEXT, not numbers, to preserve things like leading
# zeros ("011336") and trailing zeros ("0.5010"). Force text storage
# by temporarily applying a text number format before writing values.
# ---------------------------------------------------------------------
$newSheet.Range("B2:B17").NumberFormat = "@"
$newSheet.Range("D2:G17").NumberFormat = "@"

# ---- Data rows ----
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = '011336'
$newSheet.Range("C2").Value = '兴全汇吉一年持有期混合A'
$newSheet.Range("D2").Value = '15.09'
$newSheet.Range("E2").Value = '39.83'
$newSheet.Range("F2").Value = '3.32'
$newSheet.Range("G2").Value = '0.5010'
$newSheet.Range("H2").Value = 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = '011128'
$newSheet.Range("C3").Value = '华安精致生活混合A'
$newSheet.Range("D3").Value = '27.22'
$newSheet.Range("E3").Value = '87.35'
$newSheet.Range("F3").Value = '1.67'
$newSheet.Range("G3").Value = '0.4546'
$newSheet.Range("H3").Value = 10
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = '011129'
$newSheet.Range("C4").Value = '华安精致生活混合C'
$newSheet.Range("D4").Value = '13.13'
$newSheet.Range("E4").Value = '87.35'
$newSheet.Range("F4").Value = '1.67'
$newSheet.Range("G4").Value = '0.2193'
$newSheet.Range("H4").Value = 10
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = '008263'
$newSheet.Range("C5").Value = '东方红品质优选两年定期开放混合'
$newSheet.Range("D5").Value = '7.24'
$newSheet.Range("E5").Value = '27.66'
$newSheet.Range("F5").Value = '1.02'
$newSheet.Range("G5").Value = '0.0738'
$newSheet.Range("H5").Value = 7
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = '005143'
$newSheet.Range("C6").Value = '中融沪港深大消费主题灵活配置混合C'
$newSheet.Range("D6").Value = '0.69'
$newSheet.Range("E6").Value = '77.26'
$newSheet.Range("F6").Value = '8.01'
$newSheet.Range("G6").Value = '0.0553'
$newSheet.Range("H6").Value = 1
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = '003243'
$newSheet.Range("C7").Value = '上投摩根中国世纪灵活配置混合人民币份额（QDII）'
$newSheet.Range("D7").Value = '1.24'
$newSheet.Range("E7").Value = '85.53'
$newSheet.Range("F7").Value = '3.30'
$newSheet.Range("G7").Value = '0.0409'
$newSheet.Range("H7").Value = 5
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = '003244'
$newSheet.Range("C8").Value = '上投摩根中国世纪灵活配置混合美元现钞（QDII）'
$newSheet.Range("D8").Value = '1.24'
$newSheet.Range("E8").Value = '85.53'
$newSheet.Range("F8").Value = '3.30'
$newSheet.Range("G8").Value = '0.0409'
$newSheet.Range("H8").Value = 5
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = '003245'
$newSheet.Range("C9").Value = '上投摩根中国世纪灵活配置混合美元现汇（QDII）'
$newSheet.Range("D9").Value = '1.24'
$newSheet.Range("E9").Value = '85.53'
$newSheet.Range("F9").Value = '3.30'
$newSheet.Range("G9").Value = '0.0409'
$newSheet.Range("H9").Value = 5
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = '501310'
$newSheet.Range("C10").Value = '华宝标普沪港深中国增强价值指数（LOF）A'
$newSheet.Range("D10").Value = '0.96'
$newSheet.Range("E10").Value = '94.81'
$newSheet.Range("F10").Value = '2.84'
$newSheet.Range("G10").Value = '0.0273'
$newSheet.Range("H10").Value = 9
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = '011337'
$newSheet.Range("C11").Value = '兴全汇吉一年持有期混合C'
$newSheet.Range("D11").Value = '0.80'
$newSheet.Range("E11").Value = '39.83'
$newSheet.Range("F11").Value = '3.32'
$newSheet.Range("G11").Value = '0.0266'
$newSheet.Range("H11").Value = 3
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = '005142'
$newSheet.Range("C12").Value = '中融沪港深大消费主题灵活配置混合A'
$newSheet.Range("D12").Value = '0.31'
$newSheet.Range("E12").Value = '77.26'
$newSheet.Range("F12").Value = '8.01'
$newSheet.Range("G12").Value = '0.0248'
$newSheet.Range("H12").Value = 1
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = '012683'
$newSheet.Range("C13").Value = '东方红安盈甄选一年持有期混合A'
$newSheet.Range("D13").Value = '2.89'
$newSheet.Range("E13").Value = '21.56'
$newSheet.Range("F13").Value = '0.73'
$newSheet.Range("G13").Value = '0.0211'
$newSheet.Range("H13").Value = 9
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = '005701'
$newSheet.Range("C14").Value = '上投摩根香港精选港股通混合A'
$newSheet.Range("D14").Value = '0.48'
$newSheet.Range("E14").Value = '89.99'
$newSheet.Range("F14").Value = '3.16'
$newSheet.Range("G14").Value = '0.0152'
$newSheet.Range("H14").Value = 6
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = '012684'
$newSheet.Range("C15").Value = '东方红安盈甄选一年持有期混合C'
$newSheet.Range("D15").Value = '1.57'
$newSheet.Range("E15").Value = '21.56'
$newSheet.Range("F15").Value = '0.73'
$newSheet.Range("G15").Value = '0.0115'
$newSheet.Range("H15").Value = 9
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = '007397'
$newSheet.Range("C16").Value = '华宝标普沪港深中国增强价值指数（LOF）C'
$newSheet.Range("D16").Value = '0.07'
$newSheet.Range("E16").Value = '94.81'
$newSheet.Range("F16").Value = '2.84'
$newSheet.Range("G16").Value = '0.0020'
$newSheet.Range("H16").Value = 9
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = '016921'
$newSheet.Range("C17").Value = '上投摩根香港精选港股通混合C'
$newSheet.Range("D17").Value = '0.02'
$newSheet.Range("E17").Value = '89.99'
$newSheet.Range("F17").Value = '3.16'
$newSheet.Range("G17").Value = '0.0006'
$newSheet.Range("H17").Value = 6

# ---------------------------------------------------------------------
# The temporary "@" text format left an explicit style on B2:B17 and
# D2:G17; strip it back off (copy format from an untouched,
# default-styled cell) so the final cells carry no explicit style --
# matching the plain data cells used elsewhere in the workbook.
# ---------------------------------------------------------------------
$blankCell = $newSheet.Range("Z100")
$blankCell.Copy()
$newSheet.Range("B2:B17").PasteSpecial(-4122)
$newSheet.Range("D2:G17").PasteSpecial(-4122)

Write-Output "2022-Q4 sheet added and 总计 sheet updated"
